# Update "想去人数" (column F) values on sheets "展览" and "全部类型"
# to reflect the refreshed data snapshot (commit: "Update gh-pages to
# output generated at 456a3b4").

$wb = $excel.ActiveWorkbook

# New F-column values, keyed by the row number on sheet "展览".
$updates = [ordered]@{
    5  = 15717
    9  = 15450
    11 = 9048
    12 = 384
    14 = 1014
    16 = 200
    18 = 201
    20 = 57
    21 = 555
    25 = 1115
    28 = 24
    35 = 257
    36 = 326
    37 = 456
    39 = 5563
    40 = 5228
}

# "全部类型" lists the same events but has a few extra rows spliced in
# earlier in the sheet, so row numbers diverge for the later events.
$rowOnAllTypesSheet = [ordered]@{
    5  = 5
    9  = 9
    11 = 11
    12 = 12
    14 = 14
    16 = 16
    18 = 18
    20 = 20
    21 = 21
    25 = 25
    28 = 28
    35 = 37
    36 = 38
    37 = 39
    39 = 41
    40 = 43
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($row in $updates.Keys) {
    $ws1.Range("F$row").Value = $updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($row in $updates.Keys) {
    $targetRow = $rowOnAllTypesSheet[$row]
    $ws4.Range("F$targetRow").Value = $updates[$row]
}
